$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6260720491409302
$ws.Range("B1").Value = 0.481922835111618
$ws.Range("C1").Value = 0.3589450716972351
$ws.Range("D1").Value = 0.3376679122447968
$ws.Range("E1").Value = 0.3547993898391724
